# Applies scheduled market-data refresh values to the Leve profit sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1847.5
$ws.Range("J70").Value = 1847.5
$ws.Range("L70").Value = 5542.5
$ws.Range("N70").Value = -6082.5
$ws.Range("H73").Value = 1847.5
$ws.Range("J73").Value = 1847.5
$ws.Range("L73").Value = 5542.5
$ws.Range("N73").Value = -7414.5
$ws.Range("H76").Value = 3204.5454
$ws.Range("I76").Value = 3200
$ws.Range("K76").Value = 3200
$ws.Range("M76").Value = -2885
$ws.Range("H79").Value = 3204.5454
$ws.Range("I79").Value = 3200
$ws.Range("K79").Value = 3200
$ws.Range("M79").Value = -2108
$ws.Range("H88").Value = 7109.25
$ws.Range("I88").Value = 6334.3335
$ws.Range("J88").Value = 7367.5557
$ws.Range("K88").Value = 6334.3335
$ws.Range("L88").Value = 7367.5557
$ws.Range("M88").Value = -5928.3335
$ws.Range("N88").Value = -8179.5557
$ws.Range("H91").Value = 7109.25
$ws.Range("I91").Value = 6334.3335
$ws.Range("J91").Value = 7367.5557
$ws.Range("K91").Value = 6334.3335
$ws.Range("L91").Value = 7367.5557
$ws.Range("M91").Value = -4930.3335
$ws.Range("N91").Value = -10175.5557
$ws.Range("H92").Value = 847.1905
$ws.Range("I92").Value = 749.4375
$ws.Range("J92").Value = 1160
$ws.Range("K92").Value = 749.4375
$ws.Range("L92").Value = 1160
$ws.Range("M92").Value = 498.5625
$ws.Range("N92").Value = -3656
$ws.Range("H127").Value = 76923650
$ws.Range("I127").Value = 83333870
$ws.Range("K127").Value = 250001610
$ws.Range("M127").Value = -249996650
$ws.Range("H132").Value = 1951.1549
$ws.Range("I132").Value = 798.9032
$ws.Range("J132").Value = 9888.888999999999
$ws.Range("K132").Value = 2396.7096
$ws.Range("L132").Value = 29666.667
$ws.Range("M132").Value = 133.2903999999999
$ws.Range("N132").Value = -34726.667
$ws.Range("H137").Value = 2191.45
$ws.Range("I137").Value = 1447.4166
$ws.Range("J137").Value = 3307.5
$ws.Range("K137").Value = 4342.2498
$ws.Range("L137").Value = 9922.5
$ws.Range("M137").Value = -1792.2498
$ws.Range("N137").Value = -15022.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4608.226
$ws.Range("I32").Value = 2562.1096
$ws.Range("J32").Value = 18187
$ws.Range("K32").Value = 2562.1096
$ws.Range("L32").Value = 18187
$ws.Range("M32").Value = -2275.1096
$ws.Range("N32").Value = -18761
$ws.Range("H104").Value = 23498.285
$ws.Range("J104").Value = 23498.285
$ws.Range("L104").Value = 23498.285
$ws.Range("N104").Value = -30486.285
$ws.Range("H132").Value = 3196.1875
$ws.Range("I132").Value = 3112.7
$ws.Range("K132").Value = 9338.099999999999
$ws.Range("M132").Value = -6808.099999999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1252.5454
$ws.Range("I20").Value = 800.2727
$ws.Range("K20").Value = 800.2727
$ws.Range("M20").Value = -553.2727
$ws.Range("H94").Value = 1185.7142
$ws.Range("I94").Value = 1165.5555
$ws.Range("J94").Value = 1222
$ws.Range("K94").Value = 1165.5555
$ws.Range("L94").Value = 1222
$ws.Range("M94").Value = -714.5554999999999
$ws.Range("N94").Value = -2124
$ws.Range("H99").Value = 1570.9
$ws.Range("I99").Value = 1329.8572
$ws.Range("J99").Value = 2133.3333
$ws.Range("K99").Value = 1329.8572
$ws.Range("L99").Value = 2133.3333
$ws.Range("M99").Value = 168.1428000000001
$ws.Range("N99").Value = -5129.3333
$ws.Range("H107").Value = 1250.1708
$ws.Range("I107").Value = 1150.4517
$ws.Range("J107").Value = 1559.3
$ws.Range("K107").Value = 1150.4517
$ws.Range("L107").Value = 1559.3
$ws.Range("M107").Value = 769.5482999999999
$ws.Range("N107").Value = -5399.3
$ws.Range("H134").Value = 7027.909
$ws.Range("I134").Value = 6751.1665
$ws.Range("J134").Value = 7360
$ws.Range("K134").Value = 20253.4995
$ws.Range("L134").Value = 22080
$ws.Range("M134").Value = -17718.4995
$ws.Range("N134").Value = -27150

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 29413598
$ws.Range("I31").Value = 71429496
$ws.Range("J31").Value = 2471.6
$ws.Range("K31").Value = 71429496
$ws.Range("L31").Value = 2471.6
$ws.Range("M31").Value = -71429201
$ws.Range("N31").Value = -3061.6
$ws.Range("H34").Value = 29413598
$ws.Range("I34").Value = 71429496
$ws.Range("J34").Value = 2471.6
$ws.Range("K34").Value = 71429496
$ws.Range("L34").Value = 2471.6
$ws.Range("M34").Value = -71429294
$ws.Range("N34").Value = -2875.6
$ws.Range("H86").Value = 11081.8
$ws.Range("I86").Value = 6492.9
$ws.Range("J86").Value = 20259.6
$ws.Range("K86").Value = 6492.9
$ws.Range("L86").Value = 20259.6
$ws.Range("M86").Value = -5369.9
$ws.Range("N86").Value = -22505.6
$ws.Range("H89").Value = 11081.8
$ws.Range("I89").Value = 6492.9
$ws.Range("J89").Value = 20259.6
$ws.Range("K89").Value = 32464.5
$ws.Range("L89").Value = 101298
$ws.Range("M89").Value = -26848.5
$ws.Range("N89").Value = -112530
$ws.Range("H105").Value = 1437.5
$ws.Range("I105").Value = 1275
$ws.Range("J105").Value = 1600
$ws.Range("K105").Value = 1275
$ws.Range("L105").Value = 1600
$ws.Range("M105").Value = 472
$ws.Range("N105").Value = -5094
$ws.Range("H132").Value = 3769.4
$ws.Range("I132").Value = 3732.6155
$ws.Range("J132").Value = 3809.25
$ws.Range("K132").Value = 11197.8465
$ws.Range("L132").Value = 11427.75
$ws.Range("M132").Value = -8667.8465
$ws.Range("N132").Value = -16487.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1329.7
$ws.Range("I98").Value = 1816.3334
$ws.Range("J98").Value = 1121.1428
$ws.Range("K98").Value = 5449.0002
$ws.Range("L98").Value = 3363.4284
$ws.Range("M98").Value = -3951.0002
$ws.Range("N98").Value = -6359.428400000001
$ws.Range("H131").Value = 846.91
$ws.Range("I131").Value = 502.33334
$ws.Range("J131").Value = 907.71765
$ws.Range("K131").Value = 1507.00002
$ws.Range("L131").Value = 2723.15295
$ws.Range("M131").Value = 3532.99998
$ws.Range("N131").Value = -12803.15295
$ws.Range("H132").Value = 3557650
$ws.Range("I132").Value = 1589185.8
$ws.Range("J132").Value = 13892086
$ws.Range("K132").Value = 14302672.2
$ws.Range("L132").Value = 125028774
$ws.Range("M132").Value = -14300142.2
$ws.Range("N132").Value = -125033834

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1956.909
$ws.Range("I102").Value = 1802.4
$ws.Range("J102").Value = 2085.6667
$ws.Range("K102").Value = 1802.4
$ws.Range("L102").Value = 2085.6667
$ws.Range("M102").Value = -180.4000000000001
$ws.Range("N102").Value = -5329.6667
$ws.Range("H132").Value = 3248
$ws.Range("I132").Value = 3114.48
$ws.Range("K132").Value = 9343.440000000001
$ws.Range("M132").Value = -6813.440000000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 19566.666
$ws.Range("I68").Value = 100000
$ws.Range("J68").Value = 3480
$ws.Range("K68").Value = 100000
$ws.Range("L68").Value = 3480
$ws.Range("M68").Value = -99251
$ws.Range("N68").Value = -4978
$ws.Range("H71").Value = 19566.666
$ws.Range("I71").Value = 100000
$ws.Range("J71").Value = 3480
$ws.Range("K71").Value = 500000
$ws.Range("L71").Value = 17400
$ws.Range("M71").Value = -496256
$ws.Range("N71").Value = -24888

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 57500
$ws.Range("J5").Value = 57500
$ws.Range("L5").Value = 57500
$ws.Range("N5").Value = -57724
$ws.Range("H62").Value = 3000
$ws.Range("J62").Value = 3000
$ws.Range("L62").Value = 3000
$ws.Range("N62").Value = -4248
$ws.Range("H63").Value = 30247
$ws.Range("J63").Value = 30247
$ws.Range("L63").Value = 30247
$ws.Range("N63").Value = -31495
$ws.Range("H64").Value = 28375
$ws.Range("J64").Value = 28375
$ws.Range("L64").Value = 28375
$ws.Range("N64").Value = -28871
$ws.Range("H65").Value = 3000
$ws.Range("J65").Value = 3000
$ws.Range("L65").Value = 15000
$ws.Range("N65").Value = -21240
$ws.Range("H66").Value = 30247
$ws.Range("J66").Value = 30247
$ws.Range("L66").Value = 90741
$ws.Range("N66").Value = -96981
$ws.Range("H67").Value = 28375
$ws.Range("J67").Value = 28375
$ws.Range("L67").Value = 28375
$ws.Range("N67").Value = -30091
$ws.Range("H68").Value = 29633.334
$ws.Range("J68").Value = 29633.334
$ws.Range("L68").Value = 29633.334
$ws.Range("N68").Value = -31255.334
$ws.Range("H70").Value = 28999.5
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 28999.5
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 28999.5
$ws.Range("N70").Value = -29629.5
$ws.Range("M70").ClearContents()
$ws.Range("H71").Value = 29633.334
$ws.Range("J71").Value = 29633.334
$ws.Range("L71").Value = 88900.00199999999
$ws.Range("N71").Value = -97012.00199999999
$ws.Range("H73").Value = 28999.5
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 28999.5
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 28999.5
$ws.Range("N73").Value = -31183.5
$ws.Range("M73").ClearContents()
$ws.Range("H113").Value = 487.14285
$ws.Range("I113").Value = 382.22223
$ws.Range("J113").Value = 676
$ws.Range("K113").Value = 1146.66669
$ws.Range("L113").Value = 2028
$ws.Range("M113").Value = 1023.33331
$ws.Range("N113").Value = -6368
$ws.Range("H132").Value = 9528779
$ws.Range("I132").Value = 11769786
$ws.Range("J132").Value = 4499.5
$ws.Range("K132").Value = 35309358
$ws.Range("L132").Value = 13498.5
$ws.Range("M132").Value = -35306828
$ws.Range("N132").Value = -18558.5
